$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Three pairs of data rows had their match-specific columns (F:V) swapped
#    while the leading Indice/pais/torneio/temporada/data_partida columns
#    (A:E) stayed put. Swap F:V between each pair using a temp buffer.
# ---------------------------------------------------------------------------
function Swap-RowData {
    param($rowA, $rowB)
    $rangeA = $ws.Range("F$($rowA):V$($rowA)")
    $rangeB = $ws.Range("F$($rowB):V$($rowB)")
    $tmp = $rangeA.Value()
    $rangeA.Value = $rangeB.Value()
    $rangeB.Value = $tmp
}

Swap-RowData 52 53
Swap-RowData 64 65
Swap-RowData 76 77

# ---------------------------------------------------------------------------
# 2) Two brand-new match rows were appended at the bottom (rows 102 & 103).
#    Clone the formatting (styles) of the last existing row, then fill in
#    the values for each column.
# ---------------------------------------------------------------------------
$ws.Range("A101:V101").Copy()
$ws.Range("A102:V102").PasteSpecial(-4122)
$ws.Range("A103:V103").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item(102,1).Value = 101
$ws.Cells.Item(102,2).Value = "croatia"
$ws.Cells.Item(102,3).Value = "prva-nl"
$ws.Cells.Item(102,4).Value = "2023-2024"
$ws.Cells.Item(102,5).Value = 45263.53472222222
$ws.Cells.Item(102,6).Value = "Orijent"
$ws.Cells.Item(102,7).Value = 0
$ws.Cells.Item(102,8).Value = "Zrinski Jurjevac"
$ws.Cells.Item(102,9).Value = 1
$ws.Cells.Item(102,10).Value = 3.02
$ws.Cells.Item(102,11).Value = "02/12/2023 01:12"
$ws.Cells.Item(102,12).Value = 4.11
$ws.Cells.Item(102,13).Value = "03/12/2023 12:40"
$ws.Cells.Item(102,14).Value = 3.13
$ws.Cells.Item(102,15).Value = "02/12/2023 01:12"
$ws.Cells.Item(102,16).Value = 3.52
$ws.Cells.Item(102,17).Value = "03/12/2023 12:40"
$ws.Cells.Item(102,18).Value = 2.19
$ws.Cells.Item(102,19).Value = "02/12/2023 01:12"
$ws.Cells.Item(102,20).Value = 1.84
$ws.Cells.Item(102,21).Value = "03/12/2023 12:40"
$ws.Cells.Item(102,22).Value = "https://www.betexplorer.com/football/croatia/prva-nl/orijent-zrinski-jurjevac/EgFSfj00/"

$ws.Cells.Item(103,1).Value = 102
$ws.Cells.Item(103,2).Value = "croatia"
$ws.Cells.Item(103,3).Value = "prva-nl"
$ws.Cells.Item(103,4).Value = "2023-2024"
$ws.Cells.Item(103,5).Value = 45263.5625
$ws.Cells.Item(103,6).Value = "Croatia Zmijavci"
$ws.Cells.Item(103,7).Value = 1
$ws.Cells.Item(103,8).Value = "Jarun"
$ws.Cells.Item(103,9).Value = 0
$ws.Cells.Item(103,10).Value = 1.87
$ws.Cells.Item(103,11).Value = "02/12/2023 01:42"
$ws.Cells.Item(103,12).Value = 1.73
$ws.Cells.Item(103,13).Value = "03/12/2023 13:23"
$ws.Cells.Item(103,14).Value = 3.45
$ws.Cells.Item(103,15).Value = "02/12/2023 01:42"
$ws.Cells.Item(103,16).Value = 3.8
$ws.Cells.Item(103,17).Value = "03/12/2023 13:23"
$ws.Cells.Item(103,18).Value = 3.54
$ws.Cells.Item(103,19).Value = "02/12/2023 01:42"
$ws.Cells.Item(103,20).Value = 4.36
$ws.Cells.Item(103,21).Value = "03/12/2023 13:23"
$ws.Cells.Item(103,22).Value = "https://www.betexplorer.com/football/croatia/prva-nl/croatia-zmijavci-jarun/Mezl2Y7Q/"
